$wb = $excel.ActiveWorkbook

$bio = $wb.Worksheets.Item("Biology")

# Replace the text answer "46 (23 pairs)" with the plain number 46.
$bio.Range("B8").Value = 46
$bio.Range("B8").HorizontalAlignment = -4131

# Update the selection left behind on the Biology sheet and make it the
# active / visible tab (mirrors the tabSelected move away from History1).
$bio.Activate()
$bio.Range("B10").Select()
